# Update cryptocurrency price/volume table (Sheet1) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.499.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.672.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5282"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.99%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2678"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.676.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.491"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5566"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8287"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.510.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.68%  "

$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.768"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.309"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1269"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.391"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("E27").Value = "  +3.13%  "

$ws.Range("E28").Value = "  +2.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06238"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.287"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.607"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.415"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.691"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.008"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6188"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.76%  "

$ws.Range("E36").Value = "  +1.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.787"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01617"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.048"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.094.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8608"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.820.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "58.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.183"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.523"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9993"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.44%  "

$ws.Range("E50").Value = "  +0.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.15%  "
